# Applies two changes described by the commit diff:
#  1. Slide 16's table switches from table style {BA6B0F3F-66C9-4FCE-AF42-61F3C64392D9}
#     ("Table_0") to the built-in style {889AEABB-346B-47A9-B1C8-7ACBADB21AE1}.
#  2. The deck's theme color scheme (ppt/theme/theme1.xml, used by the slide
#     master / all slides) changes from the "Integral" palette to the
#     standard Office palette. (fontScheme/fmtScheme are identical between
#     the Integral and Office themes in this deck, so only clrScheme differs.)

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$s = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $tableShape = $sh
    }
}
$tableShape.Table.ApplyStyle("{889AEABB-346B-47A9-B1C8-7ACBADB21AE1}")

# --- 2. Theme colors: Integral -> Office -----------------------------------
# Order matches ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
